# Refresh the cryptos price list (Price / Volume(1h) columns) with the
# latest scraped values, mirroring the GitHub Actions data-update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.000.65"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "2.609.50"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("D5").Value = "557.76"
$ws.Range("E5").Value = "  +3.95%  "

$ws.Range("D6").Value = "144.23"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +4.44%  "

$ws.Range("D9").Value = "6.82"
$ws.Range("E9").Value = "  -3.15%  "

$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("E11").Value = "  +6.09%  "

$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").Value = "3.067.61"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").Value = "58.951.23"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").Value = "21.08"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "2.616.54"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "4.47"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").Value = "338.28"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "66.28"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("E24").Value = "  +3.14%  "

$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -1.88%  "

$ws.Range("D27").Value = "7.16"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").Value = "0.0₃0765"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").Value = "6.01"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").Value = "154.61"
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").Value = "18.94"
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").Value = "3.98"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("D35").Value = "'0.910"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.09%  "

$ws.Range("E36").Value = "  +7.97%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "37.17"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("D40").Value = "'3.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").Value = "283.93"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").Value = "0.601"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "0.0539"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "4.68"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("D49").Value = "1.946.14"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("D50").Value = "117.79"
$ws.Range("E50").Value = "  +5.70%  "

$ws.Range("D51").Value = "18.01"
$ws.Range("E51").Value = "  -1.93%  "
